$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7511
$ws1.Range("F3").Value = 70
$ws1.Range("F5").Value = 240
$ws1.Range("F6").Value = 1131
$ws1.Range("F9").Value = 131
$ws1.Range("F10").Value = 33

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7511
$ws4.Range("F3").Value = 70
$ws4.Range("F5").Value = 240
$ws4.Range("F6").Value = 1131
$ws4.Range("F10").Value = 131
$ws4.Range("F11").Value = 33
